# Add an ORGANIZATION column to the "CodeSchemes" sheet, right after the
# CODEVALUE column (i.e. a new column B, pushing ID and everything after it
# one column to the right), populate the header + the two data rows with an
# organization identifier, and grow the "yti" defined name so it still spans
# the full header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")

# Insert a new, blank column before column B (CODEVALUE stays in A, ID and
# the rest of the header shift from B.. to C..).
$ws.Range("B1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "ORGANIZATION"

# Populate the two data rows with the organization identifier.
$ws.Range("B2").Value = "74a41211-8c99-4835-a519-7a61612b1098"
$ws.Range("B3").Value = "74a41211-8c99-4835-a519-7a61612b1098"

# The "yti" defined name pointed at CodeSchemes!$A$1:$W$2 (the header range);
# since a column was inserted inside that range it now needs to reach one
# column further, to $X$2.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "yti" -or $n.Name -eq "CodeSchemes!yti") {
        $n.RefersTo = "=CodeSchemes!`$A`$1:`$X`$2"
    }
}
